$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1636683333333333
$ws.Range("H2").Value = 0.491005
$ws.Range("I2").Value = 0.008639493057305454
$ws.Range("J2").Value = 0.008639493057305455
$ws.Range("M2").Value = 110.642708
$ws.Range("N2").Value = 331.928124
$ws.Range("O2").Value = 0.5476418925386564
$ws.Range("P2").Value = 0.5476418925386564
$ws.Range("Q2").Value = 18.10870761384666
$ws.Range("R2").Value = 162.97836852462
$ws.Range("S2").Value = 0.004731348328477341
$ws.Range("T2").Value = 0.004731348328477341
$ws.Range("G3").Value = 0.1636683333333333
$ws.Range("H3").Value = 0.491005
$ws.Range("I3").Value = 0.008639493057305454
$ws.Range("J3").Value = 0.008639493057305455
$ws.Range("O3").Value = 0.3151072754333865
$ws.Range("P3").Value = 0.3151072754333865
$ws.Range("Q3").Value = 10.41955627493611
$ws.Range("R3").Value = 93.776006474425
$ws.Range("S3").Value = 0.00272236711841318
$ws.Range("T3").Value = 0.002722367118413181
$ws.Range("G4").Value = 0.1636683333333333
$ws.Range("H4").Value = 0.491005
$ws.Range("I4").Value = 0.008639493057305454
$ws.Range("J4").Value = 0.008639493057305455
$ws.Range("M4").Value = 27.72944133333333
$ws.Range("N4").Value = 83.18832399999999
$ws.Range("O4").Value = 0.1372508320279571
$ws.Range("P4").Value = 0.1372508320279571
$ws.Range("Q4").Value = 4.538431447291111
$ws.Range("R4").Value = 40.84588302562
$ws.Range("S4").Value = 0.001185777610414933
$ws.Range("T4").Value = 0.001185777610414933
$ws.Range("I5").Value = 0.808839719627903
$ws.Range("J5").Value = 0.8088397196279031
$ws.Range("M5").Value = 110.642708
$ws.Range("N5").Value = 331.928124
$ws.Range("O5").Value = 0.5476418925386564
$ws.Range("P5").Value = 0.5476418925386564
$ws.Range("Q5").Value = 1695.358962852808
$ws.Range("R5").Value = 15258.23066567527
$ws.Range("S5").Value = 0.442954514817461
$ws.Range("T5").Value = 0.4429545148174611
$ws.Range("I6").Value = 0.808839719627903
$ws.Range("J6").Value = 0.8088397196279031
$ws.Range("O6").Value = 0.3151072754333865
$ws.Range("P6").Value = 0.3151072754333865
$ws.Range("Q6").Value = 975.4913766543368
$ws.Range("S6").Value = 0.2548712803142528
$ws.Range("T6").Value = 0.2548712803142528
$ws.Range("I7").Value = 0.808839719627903
$ws.Range("J7").Value = 0.8088397196279031
$ws.Range("M7").Value = 27.72944133333333
$ws.Range("N7").Value = 83.18832399999999
$ws.Range("O7").Value = 0.1372508320279571
$ws.Range("P7").Value = 0.1372508320279571
$ws.Range("Q7").Value = 424.8934046278746
$ws.Range("R7").Value = 3824.040641650872
$ws.Range("S7").Value = 0.1110139244961893
$ws.Range("T7").Value = 0.1110139244961893
$ws.Range("G8").Value = 3.457711333333334
$ws.Range("H8").Value = 10.373134
$ws.Range("I8").Value = 0.1825207873147914
$ws.Range("J8").Value = 0.1825207873147914
$ws.Range("M8").Value = 110.642708
$ws.Range("N8").Value = 331.928124
$ws.Range("O8").Value = 0.5476418925386564
$ws.Range("P8").Value = 0.5476418925386564
$ws.Range("Q8").Value = 382.5705454022906
$ws.Range("R8").Value = 3443.134908620616
$ws.Range("S8").Value = 0.09995602939271796
$ws.Range("T8").Value = 0.09995602939271796
$ws.Range("G9").Value = 3.457711333333334
$ws.Range("H9").Value = 10.373134
$ws.Range("I9").Value = 0.1825207873147914
$ws.Range("J9").Value = 0.1825207873147914
$ws.Range("O9").Value = 0.3151072754333865
$ws.Range("P9").Value = 0.3151072754333865
$ws.Range("Q9").Value = 220.1269914979545
$ws.Range("R9").Value = 1981.14292348159
$ws.Range("S9").Value = 0.05751362800072055
$ws.Range("T9").Value = 0.05751362800072055
$ws.Range("G10").Value = 3.457711333333334
$ws.Range("H10").Value = 10.373134
$ws.Range("I10").Value = 0.1825207873147914
$ws.Range("J10").Value = 0.1825207873147914
$ws.Range("M10").Value = 27.72944133333333
$ws.Range("N10").Value = 83.18832399999999
$ws.Range("O10").Value = 0.1372508320279571
$ws.Range("P10").Value = 0.1372508320279571
$ws.Range("Q10").Value = 95.88040356526844
$ws.Range("R10").Value = 862.923632087416
$ws.Range("S10").Value = 0.02505112992135293
$ws.Range("T10").Value = 0.02505112992135293
